$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$i0 = @(8,7,7,7,7,9,5,9,8,7,8,3,9,7,8,7,8,7,6,5,7,9,6,6,7,9,8,9,9,5,9,8,8)
$if = @(8,8,8,8,7,9,6,9,8,8,8,4,9,7,8,7,8,7,6,5,7,9,6,7,8,9,8,9,9,5,9,8,8)

for ($r = 0; $r -lt 33; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
